$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that used to sit after "Git"
#    in the skills bullet list.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Italicize the school name paragraph ("Indiana University South Bend").
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`n") -eq "Indiana University South Bend") {
        $p.Range.Font.Italic = $true
        $p.Range.Font.ItalicBi = $true
        break
    }
}

# ------------------------------------------------------------------
# 3. Split the degree/date line into three runs so the attendance
#    years ("May 2015 - ") can be inserted before the graduation date,
#    and re-create the "_GoBack" bookmark right before "May 2019".
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`n") -eq "Bachelor of Science in Computer Science | May 2019") {
        $start = $p.Range.Start
        $prefix = "Bachelor of Science in Computer Science | "
        $insertPos = $start + $prefix.Length

        # Insert the new "May 2015 - " text right after the "| ".
        $insertionPoint = $d.Range($insertPos, $insertPos)
        $insertionPoint.InsertBefore("May 2015 - ")

        # Force a run boundary between "...Science | " and "May 2015 - "
        # by briefly adding and removing a bookmark at that position;
        # otherwise the two adjoining text inserts collapse into a
        # single <w:r>.
        $splitPoint = $d.Range($insertPos, $insertPos)
        $d.Bookmarks.Add("ZZZTempRunSplit", $splitPoint)
        $d.Bookmarks.Item("ZZZTempRunSplit").Delete()

        # Re-insert the "_GoBack" bookmark immediately before "May 2019".
        $bmPos = $insertPos + "May 2015 - ".Length
        $bmRange = $d.Range($bmPos, $bmPos)
        $d.Bookmarks.Add("_GoBack", $bmRange)
        break
    }
}
